$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '59.246.58'
$ws.Range('E2').Value = '  +3.06%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.111.83'
$ws.Range('E3').Value = '  +1.14%  '

$ws.Range('E4').Value = '  -0.02%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '525.02'
$ws.Range('E5').Value = '  +1.89%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '144.92'
$ws.Range('E6').Value = '  +1.46%  '

$ws.Range('E7').Value = '  -0.02%  '

$ws.Range('E8').Value = '  +1.13%  '

$ws.Range('E9').Value = '  +1.78%  '

$ws.Range('E10').Value = '  +1.54%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.386'
$ws.Range('E11').Value = '  +3.61%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '3.645.21'
$ws.Range('E12').Value = '  +0.99%  '

$ws.Range('E13').Value = '  +1.20%  '

$ws.Range('E14').Value = '  +5.09%  '

$ws.Range('E15').Value = '  +1.73%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '59.190.33'
$ws.Range('E16').Value = '  +2.76%  '

$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.108.21'
$ws.Range('E17').Value = '  +0.85%  '

$ws.Range('B18').Value = 'Polkadot'
$ws.Range('C18').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '6.22'
$ws.Range('E18').Value = '  +1.87%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.13'
$ws.Range('E19').Value = '  +0.85%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '8.23'
$ws.Range('E20').Value = '  +0.66%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '345.40'
$ws.Range('E21').Value = '  +1.93%  '

$ws.Range('E22').Value = '  -0.18%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.511'
$ws.Range('E23').Value = '  +2.13%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '66.10'
$ws.Range('E24').Value = '  +0.81%  '

$ws.Range('E25').Value = '  -0.78%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.00'
$ws.Range('E26').Value = '  -0.07%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.0₃0943'
$ws.Range('E27').Value = '  -0.08%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.87'
$ws.Range('E28').Value = '  +6.53%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.32'
$ws.Range('E29').Value = '  +2.73%  '

$ws.Range('E30').Value = '  +2.15%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.23'
$ws.Range('E31').Value = '  +3.57%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '21.15'
$ws.Range('E32').Value = '  +1.55%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '155.41'
$ws.Range('E33').Value = '  +0.71%  '

$ws.Range('E34').Value = '  +2.49%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '6.23'
$ws.Range('E35').Value = '  +5.56%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '27.01'
$ws.Range('E36').Value = '  +3.59%  '

$ws.Range('E37').Value = '  +5.23%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0693'
$ws.Range('E38').Value = '  +2.15%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.98'
$ws.Range('E39').Value = '  +3.25%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.153.74'
$ws.Range('E40').Value = '  +1.17%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '36.80'
$ws.Range('E41').Value = '  -0.59%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.669'
$ws.Range('E42').Value = '  -0.46%  '

$ws.Range('E43').Value = '  +0.00%  '

$ws.Range('E44').Value = '  +5.72%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.299.99'
$ws.Range('E45').Value = '  +1.71%  '

$ws.Range('E46').Value = '  +3.14%  '

$ws.Range('B47').Value = 'ONDO'
$ws.Range('C47').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.981'
$ws.Range('E47').Value = '  +2.48%  '

$ws.Range('B48').Value = 'InjectiveProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '20.98'
$ws.Range('E48').Value = '  +2.93%  '

$ws.Range('E49').Value = '  +3.37%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.763'
$ws.Range('E50').Value = '  +11.02%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '263.02'
$ws.Range('E51').Value = '  +11.84%  '

Write-Output "done"